$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Histórias")

# Insert a new row at position 7 for the "2.1 - Desenvolver api" task.
# This copies formatting down from row 6 and shifts rows 7-17 to 8-18.
$ws.Rows.Item(7).Insert()

# --- Fill the new row 7 ---
# Column A holds an id that looks numeric ("2.1"); force it to stay text
# like the other id cells (1.1, 1.2, ...) by formatting as Text before
# assignment, then restoring the original General/id formatting via a
# formats-only paste from a sibling id cell (avoids corrupting styles.xml).
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "2.1"
$ws.Range("A3").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B7").Value = "Desenvolver api"
$ws.Range("C7").Value = "Willians"
$ws.Range("D7").Value = "Em progresso"
$ws.Range("E7").Value = 44700
$ws.Range("F7").Value = 44700
$ws.Range("G7").Value = "Review"

# --- Update existing rows 3-6 ---
$ws.Range("D3").Value = "Em progresso"
$ws.Range("G3").Value = "Homologação"

$ws.Range("C4").Value = "Willians"
$ws.Range("D4").Value = "Em progresso"
$ws.Range("E4").Value = 44696
$ws.Range("F4").Value = 44696
$ws.Range("G4").Value = "Homologação"

$ws.Range("C5").Value = "Willians"
$ws.Range("D5").Value = "Em progresso"
$ws.Range("E5").Value = 44700
$ws.Range("F5").Value = 44700
$ws.Range("G5").Value = "Review"

$ws.Range("C6").Value = "Willians"
$ws.Range("E6").Value = 44700
$ws.Range("F6").Value = 44700
$ws.Range("G6").Value = "Desenvolvimento"

# The active selection moved to D8 (Status of "Desenvolver Cabeçalho e
# Rodapé", the row pushed down from 7 to 8) after the edits.
$ws.Range("D8").Select()
